$d = $word.ActiveDocument

$d.Content.Find.Execute("377÷9=41, 8", $true, $false, $false, $false, $false, $true, 1, $false, "237÷2=118, 1", 2) | Out-Null
$d.Content.Find.Execute("632÷5=126, 2", $true, $false, $false, $false, $false, $true, 1, $false, "249÷2=124, 1", 2) | Out-Null
$d.Content.Find.Execute("786÷9=87, 3", $true, $false, $false, $false, $false, $true, 1, $false, "679÷3=226, 1", 2) | Out-Null
$d.Content.Find.Execute("206÷8=25, 6", $true, $false, $false, $false, $false, $true, 1, $false, "730÷2=365, 0", 2) | Out-Null
$d.Content.Find.Execute("804÷3=268, 0", $true, $false, $false, $false, $false, $true, 1, $false, "893÷4=223, 1", 2) | Out-Null
$d.Content.Find.Execute("494÷9=54, 8", $true, $false, $false, $false, $false, $true, 1, $false, "252÷5=50, 2", 2) | Out-Null
$d.Content.Find.Execute("921÷8=115, 1", $true, $false, $false, $false, $false, $true, 1, $false, "235÷8=29, 3", 2) | Out-Null
$d.Content.Find.Execute("135÷8=16, 7", $true, $false, $false, $false, $false, $true, 1, $false, "995÷5=199, 0", 2) | Out-Null
$d.Content.Find.Execute("981÷9=109, 0", $true, $false, $false, $false, $false, $true, 1, $false, "965÷9=107, 2", 2) | Out-Null
$d.Content.Find.Execute("205÷8=25, 5", $true, $false, $false, $false, $false, $true, 1, $false, "268÷7=38, 2", 2) | Out-Null
$d.Content.Find.Execute("148÷7=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "335÷4=83, 3", 2) | Out-Null
$d.Content.Find.Execute("655÷5=131, 0", $true, $false, $false, $false, $false, $true, 1, $false, "264÷6=44, 0", 2) | Out-Null
$d.Content.Find.Execute("907÷7=129, 4", $true, $false, $false, $false, $false, $true, 1, $false, "719÷8=89, 7", 2) | Out-Null
$d.Content.Find.Execute("936÷2=468, 0", $true, $false, $false, $false, $false, $true, 1, $false, "428÷9=47, 5", 2) | Out-Null
$d.Content.Find.Execute("661÷6=110, 1", $true, $false, $false, $false, $false, $true, 1, $false, "601÷2=300, 1", 2) | Out-Null
$d.Content.Find.Execute("931÷7=133, 0", $true, $false, $false, $false, $false, $true, 1, $false, "583÷3=194, 1", 2) | Out-Null
$d.Content.Find.Execute("984÷5=196, 4", $true, $false, $false, $false, $false, $true, 1, $false, "988÷2=494, 0", 2) | Out-Null
$d.Content.Find.Execute("116÷7=16, 4", $true, $false, $false, $false, $false, $true, 1, $false, "243÷9=27, 0", 2) | Out-Null
$d.Content.Find.Execute("841÷2=420, 1", $true, $false, $false, $false, $false, $true, 1, $false, "492÷4=123, 0", 2) | Out-Null
$d.Content.Find.Execute("609÷9=67, 6", $true, $false, $false, $false, $false, $true, 1, $false, "290÷8=36, 2", 2) | Out-Null
$d.Content.Find.Execute("386÷4=96, 2", $true, $false, $false, $false, $false, $true, 1, $false, "594÷2=297, 0", 2) | Out-Null
$d.Content.Find.Execute("791÷9=87, 8", $true, $false, $false, $false, $false, $true, 1, $false, "446÷5=89, 1", 2) | Out-Null
$d.Content.Find.Execute("723÷3=241, 0", $true, $false, $false, $false, $false, $true, 1, $false, "869÷8=108, 5", 2) | Out-Null
$d.Content.Find.Execute("654÷6=109, 0", $true, $false, $false, $false, $false, $true, 1, $false, "159÷9=17, 6", 2) | Out-Null
$d.Content.Find.Execute("947÷6=157, 5", $true, $false, $false, $false, $false, $true, 1, $false, "948÷5=189, 3", 2) | Out-Null
